# Helper: write a value into a range while forcing it to be stored as TEXT
# (never auto-coerced into a number), mirroring how the source workbook
# stores codes like "007130" and decimal-looking figures like "40.99" as
# plain strings instead of numeric cells. The leading apostrophe forces
# Excel's text-entry semantics, then we reset the cell style back to
# "Normal" so no stray number-format is left behind on the cell.
function Set-TextCell($range, $value) {
    $range.Value = "'" + $value
    $range.Style = "Normal"
}

$wb = $excel.ActiveWorkbook

$q4Sheet = $wb.Worksheets.Item("2021-Q4")

# ------------------------------------------------------------------
# 1. Create the new "2022-Q1" sheet (a fund-holdings detail sheet),
#    positioned immediately before "总计", by duplicating the most
#    recent existing sheet of that kind ("2021-Q4") so it inherits the
#    same sheetPr/pageMargins/column styling.
#    NOTE: after Copy(Before), the reference passed as "Before" is
#    rebound to the freshly created sheet, so we grab it straight away.
# ------------------------------------------------------------------
$newSheet = $wb.Worksheets.Item("总计")
$q4Sheet.Copy($newSheet)
$newSheet.Name = "2022-Q1"

# The copied sheet has 5 data rows (rows 2-6); we need 7 data rows
# (rows 2-8), so clone row 6's formatting down into rows 7 and 8.
$newSheet.Range("A6:H6").Copy($newSheet.Range("A7:H7"))
$newSheet.Range("A6:H6").Copy($newSheet.Range("A8:H8"))

# Header row
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Fund holdings data rows
$fundRows = @(
    @(0, "007130", "中庚小盘价值股票",           "40.99", "93.10", "5.18", "2.1233", 3),
    @(1, "007497", "中庚价值灵动灵活配置混合",     "24.35", "89.42", "5.14", "1.2516", 2),
    @(2, "011174", "中庚价值品质一年持有期混合",   "63.01", "93.47", "1.86", "1.1720", 10),
    @(3, "006551", "中庚价值领航混合",             "36.49", "94.18", "2.05", "0.7480", 10),
    @(4, "002938", "中银证券健康产业灵活配置混合", "1.82",  "59.84", "4.11", "0.0748", 6),
    @(5, "004913", "中银证券聚瑞混合A",            "0.10",  "32.71", "3.44", "0.0034", 2),
    @(6, "004914", "中银证券聚瑞混合C",            "0.02",  "32.71", "3.44", "0.0007", 2)
)

$r = 2
foreach ($row in $fundRows) {
    $newSheet.Range("A$r").Value = $row[0]
    Set-TextCell $newSheet.Range("B$r") $row[1]
    Set-TextCell $newSheet.Range("C$r") $row[2]
    Set-TextCell $newSheet.Range("D$r") $row[3]
    Set-TextCell $newSheet.Range("E$r") $row[4]
    Set-TextCell $newSheet.Range("F$r") $row[5]
    Set-TextCell $newSheet.Range("G$r") $row[6]
    $newSheet.Range("H$r").Value = $row[7]
    $r = $r + 1
}

# ------------------------------------------------------------------
# 2. Update the "总计" (summary) sheet: insert a new first data row
#    for "2022-Q1" and push the existing quarters down by one row.
# ------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$totalSheet.Range("A6:D6").Copy($totalSheet.Range("A7:D7"))
$totalSheet.Range("A5:D5").Copy($totalSheet.Range("A6:D6"))
$totalSheet.Range("A4:D4").Copy($totalSheet.Range("A5:D5"))
$totalSheet.Range("A3:D3").Copy($totalSheet.Range("A4:D4"))
$totalSheet.Range("A2:D2").Copy($totalSheet.Range("A3:D3"))

$totalSheet.Range("A2").Value = 0
Set-TextCell $totalSheet.Range("B2") "2022-Q1"
$totalSheet.Range("C2").Value = 7
$totalSheet.Range("D2").Value = 5.37

$totalSheet.Range("A3").Value = 1
$totalSheet.Range("A4").Value = 2
$totalSheet.Range("A5").Value = 3
$totalSheet.Range("A6").Value = 4
$totalSheet.Range("A7").Value = 5
